# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ,
# LevePriceNQ / LevePriceHQ and the derived LeveProfitNQ / LeveProfitHQ columns
# (H:N) for the leve rows whose market data changed, per the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 440  # H2: 442.5 -> 440
$ws.Cells.Item(2, 9).Value = 128.33333  # I2: 131.66667 -> 128.33333
$ws.Cells.Item(2, 11).Value = 128.33333  # K2: 131.66667 -> 128.33333
$ws.Cells.Item(2, 13).Value = -15.33332999999999  # M2: -18.66667000000001 -> -15.33332999999999
$ws.Cells.Item(9, 8).Value = 12817  # H9: 11505.667 -> 12817
$ws.Cells.Item(9, 9).Value = 16903  # I9: 16905.5 -> 16903
$ws.Cells.Item(9, 10).Value = 559  # J9: 706 -> 559
$ws.Cells.Item(9, 11).Value = 16903  # K9: 16905.5 -> 16903
$ws.Cells.Item(9, 12).Value = 559  # L9: 706 -> 559
$ws.Cells.Item(9, 13).Value = -16734  # M9: -16736.5 -> -16734
$ws.Cells.Item(9, 14).Value = -897  # N9: -1044 -> -897
$ws.Cells.Item(10, 8).Value = 800  # H10: 700 -> 800
$ws.Cells.Item(10, 9).Value = 0  # I10: 600 -> 0
$ws.Cells.Item(10, 11).Value = 0  # K10: 600 -> 0
$ws.Cells.Item(10, 13).ClearContents()  # M10: remove (was -307)
$ws.Cells.Item(12, 8).Value = 25174.5  # H12: 25224.25 -> 25174.5
$ws.Cells.Item(12, 9).Value = 25174.5  # I12: 25224.25 -> 25174.5
$ws.Cells.Item(12, 11).Value = 25174.5  # K12: 25224.25 -> 25174.5
$ws.Cells.Item(12, 13).Value = -25004.5  # M12: -25054.25 -> -25004.5
$ws.Cells.Item(20, 8).Value = 1239  # H20: 1024.75 -> 1239
$ws.Cells.Item(20, 9).Value = 1239  # I20: 1024.75 -> 1239
$ws.Cells.Item(20, 11).Value = 1239  # K20: 1024.75 -> 1239
$ws.Cells.Item(20, 13).Value = -1009  # M20: -794.75 -> -1009
$ws.Cells.Item(21, 8).Value = 8000  # H21: 10000 -> 8000
$ws.Cells.Item(21, 9).Value = 0  # I21: 10000 -> 0
$ws.Cells.Item(21, 10).Value = 8000  # J21: 0 -> 8000
$ws.Cells.Item(21, 11).Value = 0  # K21: 10000 -> 0
$ws.Cells.Item(21, 12).Value = 8000  # L21: 0 -> 8000
$ws.Cells.Item(21, 13).ClearContents()  # M21: remove (was -9532)
$ws.Cells.Item(21, 14).Value = -8936  # N21: None -> -8936
$ws.Cells.Item(23, 8).Value = 8000  # H23: 10000 -> 8000
$ws.Cells.Item(23, 9).Value = 0  # I23: 10000 -> 0
$ws.Cells.Item(23, 10).Value = 8000  # J23: 0 -> 8000
$ws.Cells.Item(23, 11).Value = 0  # K23: 10000 -> 0
$ws.Cells.Item(23, 12).Value = 8000  # L23: 0 -> 8000
$ws.Cells.Item(23, 13).ClearContents()  # M23: remove (was -9766)
$ws.Cells.Item(23, 14).Value = -8468  # N23: None -> -8468
$ws.Cells.Item(35, 8).Value = 1239  # H35: 1024.75 -> 1239
$ws.Cells.Item(35, 9).Value = 1239  # I35: 1024.75 -> 1239
$ws.Cells.Item(35, 11).Value = 1239  # K35: 1024.75 -> 1239
$ws.Cells.Item(35, 13).Value = -860  # M35: -645.75 -> -860
$ws.Cells.Item(58, 8).Value = 3148.7  # H58: 4073.2856 -> 3148.7
$ws.Cells.Item(58, 9).Value = 673  # I58: 753.5 -> 673
$ws.Cells.Item(58, 10).Value = 6862.25  # J58: 8499.666999999999 -> 6862.25
$ws.Cells.Item(58, 11).Value = 2019  # K58: 2260.5 -> 2019
$ws.Cells.Item(58, 12).Value = 20586.75  # L58: 25499.001 -> 20586.75
$ws.Cells.Item(58, 13).Value = -1869  # M58: -2110.5 -> -1869
$ws.Cells.Item(58, 14).Value = -20886.75  # N58: -25799.001 -> -20886.75
$ws.Cells.Item(69, 8).Value = 21733.2  # H69: 22071.285 -> 21733.2
$ws.Cells.Item(69, 10).Value = 21733.2  # J69: 22071.285 -> 21733.2
$ws.Cells.Item(69, 12).Value = 65199.60000000001  # L69: 66213.855 -> 65199.60000000001
$ws.Cells.Item(69, 14).Value = -66947.60000000001  # N69: -67961.855 -> -66947.60000000001
$ws.Cells.Item(72, 8).Value = 21733.2  # H72: 22071.285 -> 21733.2
$ws.Cells.Item(72, 10).Value = 21733.2  # J72: 22071.285 -> 21733.2
$ws.Cells.Item(72, 12).Value = 195598.8  # L72: 198641.565 -> 195598.8
$ws.Cells.Item(72, 14).Value = -204334.8  # N72: -207377.565 -> -204334.8
$ws.Cells.Item(110, 8).Value = 70499  # H110: 70998.5 -> 70499
$ws.Cells.Item(110, 10).Value = 70499  # J110: 70998.5 -> 70499
$ws.Cells.Item(110, 12).Value = 70499  # L110: 70998.5 -> 70499
$ws.Cells.Item(110, 14).Value = -78679  # N110: -79178.5 -> -78679
$ws.Cells.Item(112, 8).Value = 107427.84  # H112: 89248.914 -> 107427.84
$ws.Cells.Item(112, 10).Value = 73980.78999999999  # J112: 58184.832 -> 73980.78999999999
$ws.Cells.Item(112, 12).Value = 221942.37  # L112: 174554.496 -> 221942.37
$ws.Cells.Item(112, 14).Value = -224158.37  # N112: -176770.496 -> -224158.37
$ws.Cells.Item(141, 8).Value = 2394.0833  # H141: 1898.0555 -> 2394.0833
$ws.Cells.Item(141, 9).Value = 2394.0833  # I141: 2092.4666 -> 2394.0833
$ws.Cells.Item(141, 10).Value = 0  # J141: 926 -> 0
$ws.Cells.Item(141, 11).Value = 7182.249899999999  # K141: 6277.399800000001 -> 7182.249899999999
$ws.Cells.Item(141, 12).Value = 0  # L141: 2778 -> 0
$ws.Cells.Item(141, 13).Value = -2002.249899999999  # M141: -1097.399800000001 -> -2002.249899999999
$ws.Cells.Item(141, 14).ClearContents()  # N141: remove (was -13138)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(109, 8).Value = 0  # H109: 37700 -> 0
$ws.Cells.Item(109, 10).Value = 0  # J109: 37700 -> 0
$ws.Cells.Item(109, 12).Value = 0  # L109: 37700 -> 0
$ws.Cells.Item(109, 14).ClearContents()  # N109: remove (was -40474)
$ws.Cells.Item(110, 8).Value = 143688.72  # H110: 143690.58 -> 143688.72
$ws.Cells.Item(110, 9).Value = 167487  # I110: 200807 -> 167487
$ws.Cells.Item(110, 10).Value = 899  # J110: 899.5 -> 899
$ws.Cells.Item(110, 11).Value = 167487  # K110: 200807 -> 167487
$ws.Cells.Item(110, 12).Value = 899  # L110: 899.5 -> 899
$ws.Cells.Item(110, 13).Value = -165442  # M110: -198762 -> -165442
$ws.Cells.Item(110, 14).Value = -4989  # N110: -4989.5 -> -4989
$ws.Cells.Item(122, 8).Value = 7083.4614  # H122: 7834.909 -> 7083.4614
$ws.Cells.Item(122, 9).Value = 6736.5  # I122: 7932.1665 -> 6736.5
$ws.Cells.Item(122, 10).Value = 7638.6  # J122: 7718.2 -> 7638.6
$ws.Cells.Item(122, 11).Value = 20209.5  # K122: 23796.4995 -> 20209.5
$ws.Cells.Item(122, 12).Value = 22915.8  # L122: 23154.6 -> 22915.8
$ws.Cells.Item(122, 13).Value = -17759.5  # M122: -21346.4995 -> -17759.5
$ws.Cells.Item(122, 14).Value = -27815.8  # N122: -28054.6 -> -27815.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(32, 8).Value = 21500  # H32: 0 -> 21500
$ws.Cells.Item(32, 10).Value = 21500  # J32: 0 -> 21500
$ws.Cells.Item(32, 12).Value = 21500  # L32: 0 -> 21500
$ws.Cells.Item(32, 14).Value = -22268  # N32: None -> -22268

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 7380.125  # H7: 6982 -> 7380.125
$ws.Cells.Item(7, 9).Value = 10687.363  # I7: 10737.363 -> 10687.363
$ws.Cells.Item(7, 10).Value = 104.2  # J7: 97.166664 -> 104.2
$ws.Cells.Item(7, 11).Value = 10687.363  # K7: 10737.363 -> 10687.363
$ws.Cells.Item(7, 12).Value = 104.2  # L7: 97.166664 -> 104.2
$ws.Cells.Item(7, 13).Value = -10574.363  # M7: -10624.363 -> -10574.363
$ws.Cells.Item(7, 14).Value = -330.2  # N7: -323.166664 -> -330.2
$ws.Cells.Item(22, 8).Value = 6589.647  # H22: 5676.95 -> 6589.647
$ws.Cells.Item(22, 9).Value = 7808.5713  # I22: 7306.3335 -> 7808.5713
$ws.Cells.Item(22, 10).Value = 901.3333  # J22: 788.8 -> 901.3333
$ws.Cells.Item(22, 11).Value = 7808.5713  # K22: 7306.3335 -> 7808.5713
$ws.Cells.Item(22, 12).Value = 901.3333  # L22: 788.8 -> 901.3333
$ws.Cells.Item(22, 13).Value = -7458.5713  # M22: -6956.3335 -> -7458.5713
$ws.Cells.Item(22, 14).Value = -1601.3333  # N22: -1488.8 -> -1601.3333
$ws.Cells.Item(31, 8).Value = 3817.9143  # H31: 4028.0908 -> 3817.9143
$ws.Cells.Item(31, 9).Value = 2647.3076  # I31: 2838.75 -> 2647.3076
$ws.Cells.Item(31, 11).Value = 2647.3076  # K31: 2838.75 -> 2647.3076
$ws.Cells.Item(31, 13).Value = -2352.3076  # M31: -2543.75 -> -2352.3076
$ws.Cells.Item(34, 8).Value = 3817.9143  # H34: 4028.0908 -> 3817.9143
$ws.Cells.Item(34, 9).Value = 2647.3076  # I34: 2838.75 -> 2647.3076
$ws.Cells.Item(34, 11).Value = 2647.3076  # K34: 2838.75 -> 2647.3076
$ws.Cells.Item(34, 13).Value = -2445.3076  # M34: -2636.75 -> -2445.3076
$ws.Cells.Item(58, 8).Value = 15631221  # H58: 17863854 -> 15631221
$ws.Cells.Item(58, 10).Value = 2802.2307  # J58: 2804.111 -> 2802.2307
$ws.Cells.Item(58, 12).Value = 2802.2307  # L58: 2804.111 -> 2802.2307
$ws.Cells.Item(58, 14).Value = -3208.2307  # N58: -3210.111 -> -3208.2307
$ws.Cells.Item(136, 8).Value = 15631221  # H136: 17863854 -> 15631221
$ws.Cells.Item(136, 10).Value = 2802.2307  # J136: 2804.111 -> 2802.2307
$ws.Cells.Item(136, 12).Value = 8406.6921  # L136: 8412.332999999999 -> 8406.6921
$ws.Cells.Item(136, 14).Value = -13506.6921  # N136: -13512.333 -> -13506.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 1150  # H23: 1002.1539 -> 1150
$ws.Cells.Item(23, 9).Value = 431.2  # I23: 708.6667 -> 431.2
$ws.Cells.Item(23, 10).Value = 1749  # J23: 1090.2 -> 1749
$ws.Cells.Item(23, 11).Value = 1293.6  # K23: 2126.0001 -> 1293.6
$ws.Cells.Item(23, 12).Value = 5247  # L23: 3270.6 -> 5247
$ws.Cells.Item(23, 13).Value = -1058.6  # M23: -1891.0001 -> -1058.6
$ws.Cells.Item(23, 14).Value = -5717  # N23: -3740.6 -> -5717
$ws.Cells.Item(40, 8).Value = 18.444445  # H40: 22.285715 -> 18.444445
$ws.Cells.Item(40, 10).Value = 7  # J40: 9 -> 7
$ws.Cells.Item(40, 12).Value = 28  # L40: 36 -> 28
$ws.Cells.Item(40, 14).Value = -166  # N40: -174 -> -166
$ws.Cells.Item(122, 8).Value = 1579.8  # H122: 1545 -> 1579.8
$ws.Cells.Item(122, 9).Value = 1250  # I122: 1368.6 -> 1250
$ws.Cells.Item(122, 10).Value = 1799.6666  # J122: 1671 -> 1799.6666
$ws.Cells.Item(122, 11).Value = 11250  # K122: 12317.4 -> 11250
$ws.Cells.Item(122, 12).Value = 16196.9994  # L122: 15039 -> 16196.9994
$ws.Cells.Item(122, 13).Value = -8800  # M122: -9867.4 -> -8800
$ws.Cells.Item(122, 14).Value = -21096.9994  # N122: -19939 -> -21096.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 6762.5  # H41: 7292.857 -> 6762.5
$ws.Cells.Item(41, 9).Value = 3525  # I41: 3683.3333 -> 3525
$ws.Cells.Item(41, 11).Value = 3525  # K41: 3683.3333 -> 3525
$ws.Cells.Item(41, 13).Value = -3170  # M41: -3328.3333 -> -3170
$ws.Cells.Item(126, 8).Value = 3628.4211  # H126: 3704.1052 -> 3628.4211
$ws.Cells.Item(126, 9).Value = 2402.5  # I126: 2492.375 -> 2402.5
$ws.Cells.Item(126, 11).Value = 7207.5  # K126: 7477.125 -> 7207.5
$ws.Cells.Item(126, 13).Value = -4737.5  # M126: -5007.125 -> -4737.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 6250  # H4: 11666.667 -> 6250
$ws.Cells.Item(4, 10).Value = 6250  # J4: 11666.667 -> 6250
$ws.Cells.Item(4, 12).Value = 6250  # L4: 11666.667 -> 6250
$ws.Cells.Item(4, 14).Value = -6476  # N4: -11892.667 -> -6476
$ws.Cells.Item(22, 8).Value = 3446.5833  # H22: 4042.7144 -> 3446.5833
$ws.Cells.Item(22, 9).Value = 3752.8572  # I22: 4260 -> 3752.8572
$ws.Cells.Item(22, 10).Value = 3017.8  # J22: 3499.5 -> 3017.8
$ws.Cells.Item(22, 11).Value = 3752.8572  # K22: 4260 -> 3752.8572
$ws.Cells.Item(22, 12).Value = 3017.8  # L22: 3499.5 -> 3017.8
$ws.Cells.Item(22, 13).Value = -3457.8572  # M22: -3965 -> -3457.8572
$ws.Cells.Item(22, 14).Value = -3607.8  # N22: -4089.5 -> -3607.8
$ws.Cells.Item(27, 8).Value = 3446.5833  # H27: 4042.7144 -> 3446.5833
$ws.Cells.Item(27, 9).Value = 3752.8572  # I27: 4260 -> 3752.8572
$ws.Cells.Item(27, 10).Value = 3017.8  # J27: 3499.5 -> 3017.8
$ws.Cells.Item(27, 11).Value = 3752.8572  # K27: 4260 -> 3752.8572
$ws.Cells.Item(27, 12).Value = 3017.8  # L27: 3499.5 -> 3017.8
$ws.Cells.Item(27, 13).Value = -3645.8572  # M27: -4153 -> -3645.8572
$ws.Cells.Item(27, 14).Value = -3231.8  # N27: -3713.5 -> -3231.8
$ws.Cells.Item(28, 8).Value = 6250  # H28: 11666.667 -> 6250
$ws.Cells.Item(28, 10).Value = 6250  # J28: 11666.667 -> 6250
$ws.Cells.Item(28, 12).Value = 6250  # L28: 11666.667 -> 6250
$ws.Cells.Item(28, 14).Value = -6714  # N28: -12130.667 -> -6714
$ws.Cells.Item(37, 8).Value = 6250  # H37: 11666.667 -> 6250
$ws.Cells.Item(37, 10).Value = 6250  # J37: 11666.667 -> 6250
$ws.Cells.Item(37, 12).Value = 6250  # L37: 11666.667 -> 6250
$ws.Cells.Item(37, 14).Value = -6464  # N37: -11880.667 -> -6464
$ws.Cells.Item(68, 8).Value = 2631.1667  # H68: 2647.8333 -> 2631.1667
$ws.Cells.Item(68, 9).Value = 1949.5  # I68: 1999.5 -> 1949.5
$ws.Cells.Item(68, 11).Value = 1949.5  # K68: 1999.5 -> 1949.5
$ws.Cells.Item(68, 13).Value = -1200.5  # M68: -1250.5 -> -1200.5
$ws.Cells.Item(71, 8).Value = 2631.1667  # H71: 2647.8333 -> 2631.1667
$ws.Cells.Item(71, 9).Value = 1949.5  # I71: 1999.5 -> 1949.5
$ws.Cells.Item(71, 11).Value = 9747.5  # K71: 9997.5 -> 9747.5
$ws.Cells.Item(71, 13).Value = -6003.5  # M71: -6253.5 -> -6003.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6317.5454  # H62: 6449.4 -> 6317.5454
$ws.Cells.Item(62, 10).Value = 7186.75  # J62: 7499.2856 -> 7186.75
$ws.Cells.Item(62, 12).Value = 7186.75  # L62: 7499.2856 -> 7186.75
$ws.Cells.Item(62, 14).Value = -8434.75  # N62: -8747.285599999999 -> -8434.75
$ws.Cells.Item(65, 8).Value = 6317.5454  # H65: 6449.4 -> 6317.5454
$ws.Cells.Item(65, 10).Value = 7186.75  # J65: 7499.2856 -> 7186.75
$ws.Cells.Item(65, 12).Value = 35933.75  # L65: 37496.428 -> 35933.75
$ws.Cells.Item(65, 14).Value = -42173.75  # N65: -43736.428 -> -42173.75
$ws.Cells.Item(68, 8).Value = 17499.5  # H68: 15000 -> 17499.5
$ws.Cells.Item(68, 10).Value = 17499.5  # J68: 15000 -> 17499.5
$ws.Cells.Item(68, 12).Value = 17499.5  # L68: 15000 -> 17499.5
$ws.Cells.Item(68, 14).Value = -19121.5  # N68: -16622 -> -19121.5
$ws.Cells.Item(71, 8).Value = 17499.5  # H71: 15000 -> 17499.5
$ws.Cells.Item(71, 10).Value = 17499.5  # J71: 15000 -> 17499.5
$ws.Cells.Item(71, 12).Value = 52498.5  # L71: 45000 -> 52498.5
$ws.Cells.Item(71, 14).Value = -60610.5  # N71: -53112 -> -60610.5
$ws.Cells.Item(76, 8).Value = 30000  # H76: 25000 -> 30000
$ws.Cells.Item(76, 10).Value = 30000  # J76: 25000 -> 30000
$ws.Cells.Item(76, 12).Value = 30000  # L76: 25000 -> 30000
$ws.Cells.Item(76, 14).Value = -30630  # N76: -25630 -> -30630
$ws.Cells.Item(79, 8).Value = 30000  # H79: 25000 -> 30000
$ws.Cells.Item(79, 10).Value = 30000  # J79: 25000 -> 30000
$ws.Cells.Item(79, 12).Value = 30000  # L79: 25000 -> 30000
$ws.Cells.Item(79, 14).Value = -32184  # N79: -27184 -> -32184
$ws.Cells.Item(81, 8).Value = 3859.75  # H81: 3248.25 -> 3859.75
$ws.Cells.Item(81, 9).Value = 1000  # I81: 1499 -> 1000
$ws.Cells.Item(81, 10).Value = 4813  # J81: 4997.5 -> 4813
$ws.Cells.Item(81, 11).Value = 2000  # K81: 2998 -> 2000
$ws.Cells.Item(81, 12).Value = 9626  # L81: 9995 -> 9626
$ws.Cells.Item(81, 13).Value = -939  # M81: -1937 -> -939
$ws.Cells.Item(81, 14).Value = -11748  # N81: -12117 -> -11748
$ws.Cells.Item(84, 8).Value = 3859.75  # H84: 3248.25 -> 3859.75
$ws.Cells.Item(84, 9).Value = 1000  # I84: 1499 -> 1000
$ws.Cells.Item(84, 10).Value = 4813  # J84: 4997.5 -> 4813
$ws.Cells.Item(84, 11).Value = 10000  # K84: 14990 -> 10000
$ws.Cells.Item(84, 12).Value = 48130  # L84: 49975 -> 48130
$ws.Cells.Item(84, 13).Value = -4696  # M84: -9686 -> -4696
$ws.Cells.Item(84, 14).Value = -58738  # N84: -60583 -> -58738
$ws.Cells.Item(113, 8).Value = 940.75  # H113: 964.413 -> 940.75
$ws.Cells.Item(113, 9).Value = 927.72974  # I113: 952.7714 -> 927.72974
$ws.Cells.Item(113, 10).Value = 984.5454999999999  # J113: 1001.4545 -> 984.5454999999999
$ws.Cells.Item(113, 11).Value = 2783.18922  # K113: 2858.3142 -> 2783.18922
$ws.Cells.Item(113, 12).Value = 2953.6365  # L113: 3004.3635 -> 2953.6365
$ws.Cells.Item(113, 13).Value = -613.1892200000002  # M113: -688.3141999999998 -> -613.1892200000002
$ws.Cells.Item(113, 14).Value = -7293.6365  # N113: -7344.3635 -> -7293.6365
$ws.Cells.Item(136, 8).Value = 27780356  # H136: 27780386 -> 27780356
$ws.Cells.Item(136, 9).Value = 27780356  # I136: 27780386 -> 27780356
$ws.Cells.Item(136, 11).Value = 83341068  # K136: 83341158 -> 83341068
$ws.Cells.Item(136, 13).Value = -83338518  # M136: -83338608 -> -83338518
